$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns (E:F) for the new "slrtype" fields, shifting the old
# ExcludedStudies_Excel_File_names.. ExpectedFilenames columns from E:I to G:K ---
$ws.Range("E:F").Insert()

# Match the width of the neighbouring Study_Types column (D) for the two new columns
$ws.Columns("E:F").ColumnWidth = $ws.Columns("D").ColumnWidth

# --- Populate the new / changed cells, in the same order the values were first
# introduced so that the shared-string table is built up in the expected order ---
$ws.Range("B1").Value = "Population"
$ws.Range("C1").Value = "Population_Radio_button"

$ws.Range("H3").Value = "\Testdata\Templates\ExcludedStudies\1stUpload\Excluded studies - Copy (1) - RWE.xlsx"
$ws.Range("G3").Value = "Excluded studies - Copy (1) - RWE.xlsx"

$ws.Range("F2").Value = "Clinical_radio_button"

$ws.Range("E1").Value = "slrtype"
$ws.Range("F1").Value = "slrtype_Radio_button"

$ws.Range("I3").Value = "Excluded studies - Overridedata - Copy (2).xlsx"
$ws.Range("J3").Value = "\Testdata\Templates\ExcludedStudies\Override\Excluded studies - Overridedata - Copy (2).xlsx"

$ws.Range("D2").Value = "Clinical-Interventional"
$ws.Range("D3").Value = "Clinical-RWE"

# --- Remaining cells that reuse already-known strings ---
$ws.Range("E2").Value = "Clinical"
$ws.Range("E3").Value = "Clinical"
$ws.Range("F3").Value = "Clinical_radio_button"

# --- Update the selection to match the new active cell/range ---
$ws.Range("D2:D3").Select()
